$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must remain literal text
# (matches the source data which stores everything as inlineStr).
# Force a Text number format before assigning so Excel does not
# silently coerce e.g. "1.00" -> 1.
$textForceCells = @(
    'D4'
    'D5'
    'D6'
    'D8'
    'D11'
    'D13'
    'D14'
    'D19'
    'D20'
    'D21'
    'D22'
    'D24'
    'D26'
    'D27'
    'D29'
    'D30'
    'D31'
    'D32'
    'D33'
    'D35'
    'D36'
    'D37'
    'D39'
    'D40'
    'D41'
    'D42'
    'D43'
    'D44'
    'D45'
    'D46'
    'D47'
    'D48'
    'D49'
    'D50'
    'D51'
)
foreach ($ref in $textForceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values (coin price + 1h volume change columns,
# plus the VeChain/Stellar row swap in B/C).
$values = @{
    'D2' = '63.229.63'
    'E2' = '  -0.44%  '
    'D3' = '2.643.65'
    'E3' = '  +2.62%  '
    'D4' = '1.00'
    'E4' = '  +0.13%  '
    'D5' = '604.44'
    'E5' = '  +3.08%  '
    'D6' = '144.20'
    'E6' = '  -0.20%  '
    'E7' = '  +0.11%  '
    'D8' = '0.585'
    'E8' = '  -0.31%  '
    'D9' = '2.642.09'
    'E9' = '  +2.56%  '
    'E10' = '  +1.44%  '
    'D11' = '5.64'
    'E11' = '  +0.62%  '
    'E12' = '  +0.30%  '
    'D13' = '0.363'
    'E13' = '  +3.32%  '
    'D14' = '27.28'
    'E14' = '  -0.03%  '
    'D15' = '3.124.79'
    'E15' = '  +2.81%  '
    'D16' = '63.095.62'
    'E16' = '  -0.59%  '
    'E17' = '  -0.57%  '
    'D18' = '2.657.79'
    'E18' = '  +3.25%  '
    'D19' = '11.37'
    'E19' = '  +2.39%  '
    'D20' = '4.44'
    'E20' = '  +2.98%  '
    'D21' = '342.25'
    'E21' = '  -0.20%  '
    'D22' = '6.82'
    'E22' = '  +2.85%  '
    'E23' = '  -0.17%  '
    'D24' = '67.12'
    'E24' = '  -1.63%  '
    'E25' = '  +1.39%  '
    'D26' = '1.53'
    'E26' = '  -2.78%  '
    'D27' = '8.64'
    'E27' = '  +4.95%  '
    'E28' = '  -0.08%  '
    'D29' = '545.24'
    'E29' = '  +15.83%  '
    'D30' = '1.00'
    'E30' = '  +0.34%  '
    'D31' = '7.87'
    'E31' = '  -1.00%  '
    'D32' = '2.03'
    'E32' = '  +4.22%  '
    'D33' = '1.79'
    'E33' = '  +6.70%  '
    'D34' = '0.0₃0807'
    'E34' = '  +0.74%  '
    'D35' = '172.22'
    'E35' = '  -2.46%  '
    'D36' = '5.07'
    'E36' = '  +11.50%  '
    'D37' = '0.405'
    'E37' = '  +1.74%  '
    'E38' = '  -0.11%  '
    'D39' = '19.09'
    'E39' = '  +1.20%  '
    'D40' = '1.83'
    'E40' = '  +6.49%  '
    'D41' = '172.01'
    'E41' = '  +7.69%  '
    'D42' = '1.00'
    'E42' = '  +0.06%  '
    'D43' = '3.74'
    'E43' = '  +0.82%  '
    'D44' = '22.38'
    'E44' = '  +2.99%  '
    'D45' = '0.0576'
    'E45' = '  +7.14%  '
    'D46' = '0.630'
    'E46' = '  +0.05%  '
    'B47' = 'Stellar'
    'C47' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D47' = '0.0961'
    'E47' = '  +0.14%  '
    'B48' = 'VeChain'
    'C48' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D48' = '0.0239'
    'E48' = '  +0.71%  '
    'D49' = '18.76'
    'E49' = '  +3.26%  '
    'D50' = '1.74'
    'E50' = '  +2.79%  '
    'D51' = '11.24'
    'E51' = '  -1.10%  '
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
